# Update countries & provincias Spain
# Applies the COVID data refresh (13:05 -> 13:35) to the "Pais" sheet:
#  - Updates the "Datos actualizados..." timestamp string
#  - Updates numeric stats for several countries
#  - Re-ranks Catar/Singapur and Senegal/Honduras/Bulgaria rows to keep
#    the table sorted descending by "Casos totales"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp header (row 1)
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 13:35"

# Row 11: Alemania
$ws.Range("B11").Value = 173289
$ws.Range("C11").Value = 118
$ws.Range("E11").Value = 16833
$ws.Range("G11").Value = 18
$ws.Range("H11").Value = 7756

# Row 13: Iran
$ws.Range("B13").Value = 112725
$ws.Range("C13").Value = 1958
$ws.Range("D13").Value = 89428
$ws.Range("E13").Value = 16514
$ws.Range("F13").Value = 2735
$ws.Range("G13").Value = 50
$ws.Range("H13").Value = 6783

# Rows 28-29: Catar overtakes Singapur
$ws.Range("A28").Value = "Catar"
$ws.Range("B28").Value = 26539
$ws.Range("C28").Value = 1390
$ws.Range("D28").Value = 3143
$ws.Range("E28").Value = 23382
$ws.Range("F28").Value = 72
$ws.Range("H28").Value = 14

$ws.Range("A29").Value = "Singapur"
$ws.Range("B29").Value = 25346
$ws.Range("C29").Value = 675
$ws.Range("D29").Value = 3851
$ws.Range("E29").Value = 21474
$ws.Range("F29").Value = 20
$ws.Range("H29").Value = 21

# Row 42: Filipinas
$ws.Range("F42").Value = 77

# Row 46: Dinamarca
$ws.Range("B46").Value = 10667
$ws.Range("C46").Value = 76
$ws.Range("E46").Value = 1560

# Rows 79-81: Senegal overtakes Honduras and Bulgaria
$ws.Range("A79").Value = "Senegal"
$ws.Range("B79").Value = 2105
$ws.Range("C79").Value = 110
$ws.Range("D79").Value = 782
$ws.Range("E79").Value = 1302
$ws.Range("F79").Value = 6
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 21

$ws.Range("A80").Value = "Honduras"
$ws.Range("B80").Value = 2080
$ws.Range("C80").Value = 0
$ws.Range("D80").Value = 211
$ws.Range("E80").Value = 1748
$ws.Range("F80").Value = 10
$ws.Range("G80").Value = 5
$ws.Range("H80").Value = 121

$ws.Range("A81").Value = "Bulgaria"
$ws.Range("B81").Value = 2069
$ws.Range("C81").Value = 46
$ws.Range("D81").Value = 499
$ws.Range("E81").Value = 1474
$ws.Range("F81").Value = 51
$ws.Range("G81").Value = 1
$ws.Range("H81").Value = 96

# Row 86: Republica de Macedonia
$ws.Range("B86").Value = 1694
$ws.Range("C86").Value = 20
$ws.Range("D86").Value = 1229
$ws.Range("E86").Value = 370
$ws.Range("G86").Value = 3
$ws.Range("H86").Value = 95

# Row 130: Estado de Palestina
$ws.Range("D130").Value = 310
$ws.Range("E130").Value = 63
